$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$end1 = $p1.Range.End
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$p2.Style = "Normal"

# Grab the formatted text (runs + leading empty run) of the existing
# bold "Play Chilli Heat Free..." paragraph near the end of the
# document so the new paragraph mirrors its run structure.
$pCount = $d.Paragraphs.Count
$boldSrc = $d.Paragraphs($pCount - 1)
$ft = $boldSrc.Range.FormattedText

$p2b = $d.Paragraphs(2)
$newR = $d.Range($end1, $p2b.Range.End)
$newR.FormattedText = $ft

# Re-point at the freshly created paragraph and change the bold run's
# text to "Meta description".
$p2c = $d.Paragraphs(2)
$textStart = $p2c.Range.Start
$textEnd = $p2c.Range.End - 1
$boldR = $d.Range($textStart, $textEnd)
$boldR.Text = "Meta description"

# Append the non-bold remainder of the meta description.
$afterBoldEnd = $boldR.End
$tailR = $d.Range($afterBoldEnd, $afterBoldEnd)
$tailR.InsertAfter(": Chilli Heat is a fun and engaging slot game based on Mexican cuisine featuring two bonus games for significant payouts and a 96.5% RTP. Play for free now!")
$tailR.Bold = 0

# ------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Chilli Heat Free..." paragraph
#    near the end of the document, and replace the text of the
#    following italic paragraph with the image-generation prompt.
# ------------------------------------------------------------------
$pCountEnd = $d.Paragraphs.Count
$dupBold = $d.Paragraphs($pCountEnd - 1)
$dupBold.Range.Delete()

$pCountEnd2 = $d.Paragraphs.Count
$lastP = $d.Paragraphs($pCountEnd2)
$lastTextStart = $lastP.Range.Start
$lastTextEnd = $lastP.Range.End - 1
$lastR = $d.Range($lastTextStart, $lastTextEnd)
$lastR.Text = "Prompt: Create an appealing feature image for Chilli Heat that is in line with the game's theme and features a happy Maya warrior with glasses. The image should be in a cartoon style. The feature image should include a fun and festive background, perhaps with a Mexican street party, and a cartoon version of the happy Maya warrior as the main focus. The warrior should be holding some chilli peppers and a tequila glass, with a big smile and his signature glasses. In the background, some of the symbols from the game can be included, such as the mariachi, the chihuahua with the tabasco sauce, and the sacks of coins. The overall feel of the image should be colorful and engaging, with a touch of humor to reflect the fun and laid-back nature of the game."

Write-Output "done"
